$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.433413505554199
$ws.Range("B1").Value = 5.007958889007568
$ws.Range("C1").Value = 6.829205989837646
$ws.Range("D1").Value = 9.128503799438477
$ws.Range("E1").Value = 4.804558277130127
